$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove duplicated teacher-list values, replacing them with a plain "-"
$ws.Range("C7").Value = "-"
$ws.Range("D18").Value = "-"
$ws.Range("D19").Value = "-"
$ws.Range("D20").Value = "-"
$ws.Range("D21").Value = "-"
